# Applies the "Update excel models and documentation" revision to tasks.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sharedStrings edit: "Analysis" -> "Backend Analysis" (TaskName for row 9) ---
$ws.Range("D9").Value = "Backend Analysis"

# --- new TaskDependencies (column F) values for a few rows ---
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = 4

# --- rows 5 & 6: TaskName/Hours/InternalID corrected (effectively swapped) ---
$ws.Range("D5").Value = "Database Setup"
$ws.Range("E5").Value = 60
$ws.Range("I5").Value = 1327

$ws.Range("D6").Value = "API Development"
$ws.Range("E6").Value = 50
$ws.Range("I6").Value = 1236

# --- EstimatedEffortHours corrections ---
$ws.Range("E9").Value = 47
$ws.Range("E11").Value = 45

# --- table still covers the full data range A1:I13 ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I13"))

# --- view state: zoom to 70% and move the active selection to G7 ---
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("G7").Select()
